$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '22.530.03'
Set-TextValue $ws 'E2' '  +0.21%  '
Set-TextValue $ws 'D3' '1.577.03'
Set-TextValue $ws 'E3' '  +0.27%  '
Set-TextValue $ws 'D4' '1.003'
Set-TextValue $ws 'E4' '  +0.10%  '
Set-TextValue $ws 'D6' '288.52'
Set-TextValue $ws 'E6' '  -0.96%  '
Set-TextValue $ws 'D7' '0.3699'
Set-TextValue $ws 'E7' '  -0.15%  '
Set-TextValue $ws 'D8' '48.64'
Set-TextValue $ws 'E8' '  -2.47%  '
Set-TextValue $ws 'E9' '  -1.00%  '
Set-TextValue $ws 'D10' '1.142'
Set-TextValue $ws 'E10' '  +0.16%  '
Set-TextValue $ws 'D11' '0.07481'
Set-TextValue $ws 'E11' '  -0.79%  '
Set-TextValue $ws 'D12' '1.003'
Set-TextValue $ws 'E12' '  +0.11%  '
Set-TextValue $ws 'D13' '20.99'
Set-TextValue $ws 'E13' '  -1.25%  '
Set-TextValue $ws 'D14' '6.003'
Set-TextValue $ws 'E14' '  -0.26%  '
Set-TextValue $ws 'D15' '6.963'
Set-TextValue $ws 'E15' '  +0.15%  '
Set-TextValue $ws 'D16' '1.576.96'
Set-TextValue $ws 'E16' '  +0.32%  '
Set-TextValue $ws 'D17' '0.00001120'
Set-TextValue $ws 'E17' '  +0.14%  '
Set-TextValue $ws 'D18' '88.73'
Set-TextValue $ws 'E18' '  -2.10%  '
Set-TextValue $ws 'D19' '0.06766'
Set-TextValue $ws 'B20' 'Dai'
Set-TextValue $ws 'C20' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D20' '1.002'
Set-TextValue $ws 'E20' '  +0.07%  '
Set-TextValue $ws 'B21' 'Uniswap'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws 'D21' '6.413'
Set-TextValue $ws 'E21' '  +1.75%  '
Set-TextValue $ws 'D22' '16.56'
Set-TextValue $ws 'E22' '  +0.72%  '
Set-TextValue $ws 'D23' '12.16'
Set-TextValue $ws 'E23' '  -0.59%  '
Set-TextValue $ws 'D24' '22.534.08'
Set-TextValue $ws 'E24' '  +0.20%  '
Set-TextValue $ws 'D25' '2.396'
Set-TextValue $ws 'E25' '  +0.98%  '
Set-TextValue $ws 'D26' '2.599'
Set-TextValue $ws 'E26' '  +0.00%  '
Set-TextValue $ws 'D27' '152.74'
Set-TextValue $ws 'E27' '  +2.34%  '
Set-TextValue $ws 'D28' '19.70'
Set-TextValue $ws 'E28' '  -1.85%  '
Set-TextValue $ws 'D29' '5.011'
Set-TextValue $ws 'D30' '124.48'
Set-TextValue $ws 'E30' '  -0.47%  '
Set-TextValue $ws 'D31' '1.757.19'
Set-TextValue $ws 'E31' '  +0.50%  '
Set-TextValue $ws 'D32' '1.065'
Set-TextValue $ws 'E32' '  -0.64%  '
Set-TextValue $ws 'D33' '6.188'
Set-TextValue $ws 'E33' '  -0.63%  '
Set-TextValue $ws 'D34' '2.007'
Set-TextValue $ws 'E34' '  -0.15%  '
Set-TextValue $ws 'D35' '9.666'
Set-TextValue $ws 'E35' '  -0.97%  '
Set-TextValue $ws 'D36' '0.08341'
Set-TextValue $ws 'E36' '  -0.05%  '
Set-TextValue $ws 'D37' '0.02460'
Set-TextValue $ws 'E37' '  -0.99%  '
Set-TextValue $ws 'D38' '0.2275'
Set-TextValue $ws 'E38' '  -1.18%  '
Set-TextValue $ws 'D39' '5.445'
Set-TextValue $ws 'E39' '  +0.17%  '
Set-TextValue $ws 'D40' '0.06388'
Set-TextValue $ws 'E40' '  -2.36%  '
Set-TextValue $ws 'D41' '1.296'
Set-TextValue $ws 'E41' '  -4.88%  '
Set-TextValue $ws 'D42' '0.6361'
Set-TextValue $ws 'E42' '  +2.29%  '
Set-TextValue $ws 'D43' '11.38'
Set-TextValue $ws 'E43' '  +0.51%  '
Set-TextValue $ws 'E44' '  +0.10%  '
Set-TextValue $ws 'D45' '14.05'
Set-TextValue $ws 'E45' '  -0.29%  '
Set-TextValue $ws 'D46' '0.6195'
Set-TextValue $ws 'E46' '  +5.86%  '
Set-TextValue $ws 'D47' '3.771'
Set-TextValue $ws 'E47' '  -0.89%  '
Set-TextValue $ws 'D48' '2.064'
Set-TextValue $ws 'E48' '  -0.35%  '
Set-TextValue $ws 'D49' '124.89'
Set-TextValue $ws 'E49' '  -3.29%  '
Set-TextValue $ws 'D50' '1.217'
Set-TextValue $ws 'E50' '  -0.46%  '
Set-TextValue $ws 'D51' '0.07274'
Set-TextValue $ws 'E51' '  -0.75%  '

$wb.Save()
Write-Output "Done applying updates."
